$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 (patient record "LEO MESSI" -> "KILLIAN MBPAPPE") ---
$ws.Range("B5").Value = "13-12-2023"
$ws.Range("C5").Value = "KILLIAN"
$ws.Range("D5").Value = "MBPAPPE"

# F5 and M5 look like numbers (leading zero / plain digits) - force text so
# Excel doesn't coerce them into numeric cells (matches original inlineStr text).
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "07894125632"
$ws.Range("F5").Style = "Normal"

$ws.Range("I5").Value = "IF-ILFOV"
$ws.Range("K5").Value = "Pensionar"

$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "256312"
$ws.Range("M5").Style = "Normal"

$ws.Range("N5").Value = "ASDSADSA`n"
$ws.Range("T5").Value = "YES"
$ws.Range("U5").Value = "EGOISM`n"
$ws.Range("V5").Value = "A SE MUTA LA REAL MADRID`n"

# Re-fit the row after writing multi-line text so no stray custom row-height
# sticks around in the saved file.
$ws.Rows(5).EntireRow.AutoFit()

# --- Remove row 6 (VASILE GHEORGHE record) entirely ---
$ws.Rows(6).Delete()

# --- Column width tweaks ---
# ColumnWidth (character units) differs from the stored XML width by a fixed
# +5/6 padding offset, so subtract that to land on the exact target width.
$ws.Columns(6).ColumnWidth = 21 - 5/6    # F: 20 -> 21
$ws.Columns(8).ColumnWidth = 20 - 5/6    # H: 24 -> 20
$ws.Columns(9).ColumnWidth = 20 - 5/6    # I: 21 -> 20
$ws.Columns(21).ColumnWidth = 55 - 5/6   # U: 62 -> 55

# --- Refresh the AutoFilter range to the new used range (A1:V5) ---
$ws.AutoFilterMode = $false
$ws.Range("A1:V5").AutoFilter() | Out-Null

# --- Fix up the _FilterDatabase defined name to match the new range ---
foreach ($n in $wb.Names) {
  if ($n.Name -eq "REGISTRU!_FilterDatabase") {
    $n.RefersTo = "='REGISTRU'!`$A`$1:`$V`$5"
  }
}
